# Update cryptocurrency price/volume figures to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.33%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.09%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.124"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.23%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08117"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.19%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.946"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.10%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.146"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.73%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9268"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.38%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1408"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.55%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1936"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.31%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09036"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.41%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03495"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.31%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09824"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.18%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001394"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.72%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006041"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.93%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.905"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.236"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.89%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.16%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.36%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.744"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.30%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.72%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04377"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.53%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001231"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.97%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004798"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.19%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.07%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004003"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.99%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02083"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.31%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05117"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.41%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007434"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.42%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009786"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.01%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1363"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.04%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002132"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.53%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009205"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.32%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006395"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.80%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.08%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-18.95%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002575"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.08%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.08%"
